$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($row in 2,3) {
    $ws.Range("D$row").Value = -0.0582
    $ws.Range("G$row").Value = -0.1265350877192983
    $ws.Range("H$row").Value = -0.1265350877192983
    $ws.Range("I$row").Value = -0.3793859649122807
    $ws.Range("J$row").Value = -0.3793859649122807
    $ws.Range("K$row").Value = -1.57
    $ws.Range("L$row").Value = -0.3442982456140351
    $ws.Range("U$row").Value = 0.011
    $ws.Range("V$row").Value = 0.002716049382716049
    $ws.Range("W$row").Value = -0.1880239520958084
    $ws.Range("X$row").Value = 0.09802955230621008
    $ws.Range("Y$row").Value = -0.2860535044020185
    $ws.Range("Z$row").Value = 0.4534606205250596
    $ws.Range("AA$row").Value = -0.1720365950676213
    $ws.Range("AB$row").Value = 0.08031034660618078
    $ws.Range("AC$row").Value = -0.2523469416738021
    $ws.Range("AD$row").Value = 1.5
    $ws.Range("AF$row").Value = 1.5
    $ws.Range("AG$row").Value = 1.489
    $ws.Range("AH$row").Value = 0.2702702702702703
    $ws.Range("AI$row").Value = 0.1704545454545454
    $ws.Range("AJ$row").Value = 0.2688210868387796
    $ws.Range("AK$row").Value = 0.1694163158493572
    $ws.Range("AL$row").Value = 0.043
    $ws.Range("AM$row").Value = 0.043
    $ws.Range("AN$row").Value = -1.209677419354839
    $ws.Range("AO$row").Value = -40.23255813953489
    $ws.Range("AP$row").Value = -1.200806451612903
    $ws.Range("AQ$row").Value = -40.23255813953489
}
